# Edit: LevantamentoRequisitos - add new RC07-RC10/RE01-RE12 requirement rows
# and revise two existing requirement descriptions (RD05/RC05 text, RC06 text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DocumentoGeral")

# --- Revise existing requirement descriptions ---

# RC05 description (row 22): "...adicionar novos items." -> "...adicionar e alterar itens."
$ws.Cells.Item(22, 3).Value = 'Os administrativos podem atualizar o sálario dos funcionários, alterar o estado dos pedidos e adicionar e alterar itens.'

# RD06 description (row 24): distribuidor -> funcionarios
$ws.Cells.Item(24, 3).Value = 'A cada pedido estão associados varios funcionários e a cada funcionario muitos pedidos.'

# --- Append new requirement rows (17-32) ---

$ws.Cells.Item(27, 2).Value = 17
$ws.Cells.Item(27, 3).Value = 'Somente Jefferson Bazos tem controle sobre a informação completa da base de dados'
$ws.Cells.Item(27, 5).Value = 'RC07'

$ws.Cells.Item(28, 2).Value = 18
$ws.Cells.Item(28, 3).Value = 'Somente Jefferson bazos pode ter acessso a informação em relação aos gestores dos armazens '
$ws.Cells.Item(28, 5).Value = 'RC08'

$ws.Cells.Item(29, 2).Value = 19
$ws.Cells.Item(29, 3).Value = 'Somente o gestor do armazém e Jefferson bazos pode acessar os dados do seu armazem (funcionarios, itens , pedidos)'
$ws.Cells.Item(29, 5).Value = 'RC09'

$ws.Cells.Item(30, 2).Value = 20
$ws.Cells.Item(30, 3).Value = 'Obter o número de clientes '
$ws.Cells.Item(30, 4).Value = 'Clientes'
$ws.Cells.Item(30, 5).Value = 'RE01'

$ws.Cells.Item(31, 2).Value = 21
$ws.Cells.Item(31, 3).Value = 'Verificar o número de funcionários '
$ws.Cells.Item(31, 4).Value = 'Funcionários'
$ws.Cells.Item(31, 5).Value = 'RE02'

$ws.Cells.Item(32, 2).Value = 22
$ws.Cells.Item(32, 3).Value = 'Saber para cada pedido quais funcionarios estão associados'
$ws.Cells.Item(32, 4).Value = 'Pedidos/Funcionários'
$ws.Cells.Item(32, 5).Value = 'RE03'

$ws.Cells.Item(33, 2).Value = 23
$ws.Cells.Item(33, 3).Value = 'Conhecer o armazém em que o pedido está relacionado'
$ws.Cells.Item(33, 4).Value = 'Armazém/Pedido'
$ws.Cells.Item(33, 5).Value = 'RE04'

$ws.Cells.Item(34, 2).Value = 24
$ws.Cells.Item(34, 3).Value = 'Identificar os clientes com maiores gastos'
$ws.Cells.Item(34, 4).Value = 'Clientes'
$ws.Cells.Item(34, 5).Value = 'RE05'

$ws.Cells.Item(35, 2).Value = 25
$ws.Cells.Item(35, 3).Value = 'Consultar os pedidos feitos por um cliente'
$ws.Cells.Item(35, 4).Value = 'Clientes/Pedidos'
$ws.Cells.Item(35, 5).Value = 'RE06'

$ws.Cells.Item(36, 2).Value = 26
$ws.Cells.Item(36, 3).Value = 'Consultar os funcionários com melhores desempenho'
$ws.Cells.Item(36, 4).Value = 'Funcionários'
$ws.Cells.Item(36, 5).Value = 'RE07'

$ws.Cells.Item(37, 2).Value = 27
$ws.Cells.Item(37, 3).Value = 'Identificar os items mais vendidos'
$ws.Cells.Item(37, 4).Value = 'Items'
$ws.Cells.Item(37, 5).Value = 'RE08'

$ws.Cells.Item(38, 2).Value = 28
$ws.Cells.Item(38, 3).Value = 'Saber para cada funcionario quais pedidos estão associados'
$ws.Cells.Item(38, 4).Value = 'Funcionários'
$ws.Cells.Item(38, 5).Value = 'RE09'

$ws.Cells.Item(39, 2).Value = 29
$ws.Cells.Item(39, 3).Value = 'Obter quais funcionários trabalham em dado armazem'
$ws.Cells.Item(39, 4).Value = 'Funcionários'
$ws.Cells.Item(39, 5).Value = 'RE10'

$ws.Cells.Item(40, 2).Value = 30
$ws.Cells.Item(40, 3).Value = 'Obter quais tipos de items são mais populares'
$ws.Cells.Item(40, 4).Value = 'Items'
$ws.Cells.Item(40, 5).Value = 'RE11'

$ws.Cells.Item(41, 2).Value = 31
$ws.Cells.Item(41, 3).Value = 'Obter uma relação entre salário e desempenho de cada funcionarios por tipo'
$ws.Cells.Item(41, 4).Value = 'Funcionários'
$ws.Cells.Item(41, 5).Value = 'RE12'

$ws.Cells.Item(42, 2).Value = 32
$ws.Cells.Item(42, 3).Value = 'Os funcionarios relacionados a um pedido tem acesso a informação desse pedido'
$ws.Cells.Item(42, 4).Value = 'Funcionários/Pedido'
$ws.Cells.Item(42, 5).Value = 'RC10'

# --- Cursor/selection bookkeeping to mirror the final workbook state ---

# Requisitos_Manipulação keeps its own last-used cell selection
$wsManip = $wb.Worksheets.Item("Requisitos_Manipulação")
$wsManip.Range("C18").Select() | Out-Null

# Restore DocumentoGeral as the active sheet, with its cursor on the last
# filled cell of the new block
$ws.Activate() | Out-Null
$ws.Range("E43").Select() | Out-Null
